$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (columns C and D got narrower) ---
$ws.Columns.Item(3).ColumnWidth = 13.75
$ws.Columns.Item(4).ColumnWidth = 13.75

# --- Header row (row 1) ---
# Row 1
$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Prognose"
$ws.Range("C1").Value = "surveys"
$ws.Range("D1").Value = "production"
$ws.Range("E1").Value = "orders"
$ws.Range("F1").Value = "turnover"
$ws.Range("G1").Value = "financial"
$ws.Range("H1").Value = "labor market"
$ws.Range("I1").Value = "prices"
$ws.Range("J1").Value = "national accounts"
$ws.Range("K1").Value = "Revision"

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-03-30"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = 0.29515230513940172
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-04-15"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = 0.29590581634224794
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0011378175204363643
$ws.Range("E3").Value = -0.00007482977947776712
$ws.Range("F3").Value = 0.000002780734069021098
$ws.Range("G3").Value = -0.000045767141353086311
$ws.Range("H3").Value = 0.00000032361616386279458
$ws.Range("I3").Value = -0.00034526102479022025
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.000078447277798043835

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2025-04-30"
$ws.Range("A4").Style = "Normal"
$ws.Range("B4").Value = 0.28563153069162722
$ws.Range("C4").Value = -0.0092187592196895981
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = -0.000010552066146064088
$ws.Range("F4").Value = -0.00000038252996974079545
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.0000067906626925608551
$ws.Range("I4").Value = -0.0014760168017231134
$ws.Range("J4").Value = 0.00042033815336391057
$ws.Range("K4").Value = 0.0000042961508513461588

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-05-15"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 0.29018361655644137
$ws.Range("C5").Value = 0.0077312126418730255
$ws.Range("D5").Value = 0.0010170571814887128
$ws.Range("E5").Value = -0.00028040171741685623
$ws.Range("F5").Value = -0.00054262425465993714
$ws.Range("G5").Value = -0.0024635687512932322
$ws.Range("H5").Value = -0.000055471590411047531
$ws.Range("I5").Value = -0.0005368735911450152
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.00031724405362149044

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-05-30"
$ws.Range("A6").Style = "Normal"
$ws.Range("B6").Value = 0.56973772635444619
$ws.Range("C6").Value = 0.28173281508188691
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.000028075922827868678
$ws.Range("F6").Value = -0.00036127241602645904
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -0.000066647633447197
$ws.Range("I6").Value = -0.0055867981475141828
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.0038079369902778715

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-06-15"
$ws.Range("A7").Style = "Normal"
$ws.Range("B7").Value = 0.53429475734826581
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = -0.038656790983625983
$ws.Range("E7").Value = -0.00031369078075605978
$ws.Range("F7").Value = -0.0035085940873540882
$ws.Range("G7").Value = 0.0055177399755155275
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0.0014233120359632674
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0.000095054834076901962

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "2025-06-30"
$ws.Range("A8").Style = "Normal"
$ws.Range("B8").Value = 0.28053051091238096
$ws.Range("C8").Value = -0.25365022355844302
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.000060934941984182012
$ws.Range("F8").Value = -0.0094807057218160521
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = -0.000079753525788405375
$ws.Range("I8").Value = 0.0091161670030319615
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0.00026933442514642802

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "2025-07-15"
$ws.Range("A9").Style = "Normal"
$ws.Range("B9").Value = -0.0064011239419219379
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = -0.06089166780733344
$ws.Range("E9").Value = -0.036499888744238179
$ws.Range("F9").Value = -0.18170754325635502
$ws.Range("G9").Value = -0.0019614745688454655
$ws.Range("H9").Value = -0.0067467602845291807
$ws.Range("I9").Value = 0.00080904066026612825
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0.00006665914673226192

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-07-30"
$ws.Range("A10").Style = "Normal"
$ws.Range("B10").Value = 0.28153244975661457
$ws.Range("C10").Value = 0.30727403193837344
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = -0.00048396533642513384
$ws.Range("F10").Value = 0.0052132085936798418
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = -0.00063712489764264607
$ws.Range("I10").Value = 0.0033691202965397357
$ws.Range("J10").Value = -0.024230238910139673
$ws.Range("K10").Value = -0.0025714579858490683

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "2025-08-15"
$ws.Range("A11").Style = "Normal"
$ws.Range("B11").Value = 0.36807994259692067
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.016123389531590371
$ws.Range("E11").Value = 0.021510082238487788
$ws.Range("F11").Value = 0.11709711786347254
$ws.Range("G11").Value = 0.0044722874112348814
$ws.Range("H11").Value = 0.00080038251862255385
$ws.Range("I11").Value = -0.015211784396545911
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = -0.058243982326556115

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "2025-08-30"
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Value = 0.3195919305058359
$ws.Range("C12").Value = -0.0205433504581066
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0.00047162040210552364
$ws.Range("F12").Value = 0.00014335484008250762
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = -0.00037814925371547705
$ws.Range("I12").Value = -0.0084832192442814095
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = -0.019698268377169315

Write-Host "applied nowcast update"